# The underlying data rows (10-32, row 19 untouched) in the "Artfynd" sheet
# were re-shuffled/re-ordered by the source system: the content that used to
# sit in one row now belongs to another row (a pure row-content permutation
# over columns A:AY). Row 19 keeps its original content.
#
# Mapping: destination row -> source row (content to copy from, read from the
# ORIGINAL/before state).
#   10 <- 26   11 <- 20   12 <- 23   13 <- 25   14 <- 27
#   15 <- 22   16 <- 11   17 <- 12   18 <- 14   20 <- 21
#   21 <- 15   22 <- 17   23 <- 28   24 <- 29   25 <- 16
#   26 <- 10   27 <- 18   28 <- 13   29 <- 30   30 <- 24
#   31 <- 32   32 <- 31

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = "A"
$lastCol = "AY"

# Columns that hold date/time-looking text (Startdatum, Starttid, Slutdatum,
# Sluttid) which Excel would otherwise auto-coerce into real date serials
# when a plain .Value assignment sees a string like "2023-08-14". Forcing
# text format on the destination before writing keeps them as plain strings,
# matching the workbook's original inlineStr encoding.
# (column indices within the A:AY range: Y=25, Z=26, AA=27, AB=28)
$dateColIndexes = @(25, 26, 27, 28)

# 1) Snapshot every source row (10-32) BEFORE any writes, since several rows
#    participate in permutation cycles and would otherwise clobber data that
#    another destination still needs to read.
$snapshot = @{}
for ($r = 10; $r -le 32; $r++) {
    $snapshot[$r] = $ws.Range("$firstCol$r`:$lastCol$r").Value()
}

$rowMap = @{
    10 = 26; 11 = 20; 12 = 23; 13 = 25; 14 = 27;
    15 = 22; 16 = 11; 17 = 12; 18 = 14;
    20 = 21; 21 = 15; 22 = 17; 23 = 28; 24 = 29; 25 = 16;
    26 = 10; 27 = 18; 28 = 13;
    29 = 30; 30 = 24;
    31 = 32; 32 = 31
}

# 2) Write the snapshotted content into its new destination row. Row 19 is
#    intentionally left untouched (it maps to itself).
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $destRange = $ws.Range("$firstCol$destRow`:$lastCol$destRow")
    foreach ($idx in $dateColIndexes) {
        $destRange.Columns.Item($idx).NumberFormat = "@"
    }
    $destRange.Value = $snapshot[$srcRow]
}
